$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B currently holds text values "s"/"a"/"s"/"s" (stored as shared strings).
# Replace them with numeric values: B1=1, B2=1, B3=1, B4=2
$ws.Range("B1").Value = 1
$ws.Range("B2").Value = 1
$ws.Range("B3").Value = 1
$ws.Range("B4").Value = 2
